# "yetki rol eklentisi ama yarım biraz bozuk ayrıca kontrollerlar yok"
#
# Insert 5 new columns (A:E) in front of the existing "db yapısı" table
# so everything currently in columns A:S shifts to F:X, then add the
# new "rol" / "ekran" / "ekran_rol" / "rol_kisi" table headers in the
# freshly inserted A1:D1 (E stays empty, matching the target layout).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the whole existing sheet 5 columns to the right (A:S -> F:X).
$ws.Range("A1:E1").EntireColumn.Insert()

# New header row cells, added in this exact order so the new shared
# strings come out as rol, ekran, ekran_rol, rol_kisi.
$ws.Range("A1").Value = "rol"
$ws.Range("C1").Value = "ekran"
$ws.Range("D1").Value = "ekran_rol"
$ws.Range("B1").Value = "rol_kisi"

# Match the author's new selection.
$ws.Range("K16").Select()
